# Bubble sheet answers: add 5 new respondent rows (139-143). Rows 141/142
# duplicate the answer grids currently on rows 136/137, row 143 duplicates
# row 138 (including its blank-but-present Name/B cell), the stray empty
# Name cell is removed from the now-interior row 138, and two brand new
# respondent rows (139, 140) are inserted ahead of the duplicated tail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C..AF in order, used for the per-row numeric answer grids below.
$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF")

function Set-RowAnswers($row, $values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# --- Row 143 inherits row 138's blank-but-present Name (B) cell: copy just
# --- that single cell over before row 138 loses it below.
$ws.Range("B138").Copy($ws.Range("B143"))

# --- Row 138 no longer carries the stray empty Name (B) cell.
$ws.Range("B138").ClearContents()

# --- Row 139: "826" / "AbnedAK"
$ws.Range("A139").Value = "'826"
$ws.Range("B139").Value = "AbnedAK"
Set-RowAnswers 139 @(0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0)

# --- Row 140: "810"
$ws.Range("A140").Value = "'810"
Set-RowAnswers 140 @(1,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0)

# --- Row 141: "8"
$ws.Range("A141").Value = "'8"
Set-RowAnswers 141 @(1,0,0,0,1,1,1,1,0,0,0,0,0,0,1,0,0,1,0,0,1,0,0,0,0,0,1,0,0,0)

# --- Row 142: "79370"
$ws.Range("A142").Value = "'79370"
Set-RowAnswers 142 @(1,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0)

# --- Row 143: "0082018"
$ws.Range("A143").Value = "'0082018"
Set-RowAnswers 143 @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0)
